$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# First draft of scenario-specific input tables: bump the per-year cost
# figures in rows 2 and 3 (columns D:AR) from 10 to 50.
$ws.Range("D2:AR3").Value = 50

# Leave the range selected with D2 as the active cell, matching the
# saved workbook's UI state.
$ws.Range("D2:AR3").Select()
